# Daily attendance processing - reorders the comma-separated "Recorded By"
# names in column G so that higher-priority actors (admin, then
# System/system) sort before the ordinary recorder e-mail addresses,
# preserving the relative order of any remaining names (stable sort).

function Test-ExactMatch($a, $b) {
    # Ordinal/case-sensitive comparison via .CompareTo (the -eq/-ceq
    # operators in this host are case-insensitive, so we can't rely on them
    # to tell "System" from "system").
    return ($a.CompareTo($b) -eq 0)
}

function Get-NameRank($name) {
    if (Test-ExactMatch $name "admin@admin.com") { return 0 }
    if (Test-ExactMatch $name "System") { return 1 }
    if (Test-ExactMatch $name "system") { return 1.5 }
    return 2
}

function Reorder-RecordedBy($value) {
    $parts = $value -split ", "
    if ($parts.Length -le 1) {
        return $value
    }

    $indexed = @()
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $part = $parts[$i]
        $rank = Get-NameRank $part
        $indexed += , @($rank, $i, $part)
    }

    # Stable sort: primary key = priority rank, secondary key = original
    # position (keeps ties in their original relative order).
    $sorted = $indexed | Sort-Object { $_[0] }, { $_[1] }

    $resultParts = @()
    foreach ($item in $sorted) {
        $resultParts += $item[2]
    }

    return ($resultParts -join ", ")
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2
    if ($current -eq $null) {
        continue
    }
    $updated = Reorder-RecordedBy $current
    if ($updated -ne $current) {
        $cell.Value = $updated
    }
}
